$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 50000176
$ws.Range("I8").Value = 52631604
$ws.Range("K8").Value = 157894812
$ws.Range("M8").Value = -157894673
$ws.Range("H19").Value = 3038.7334
$ws.Range("J19").Value = 1547
$ws.Range("L19").Value = 1547
$ws.Range("N19").Value = -1897
$ws.Range("H137").Value = 3179750.2
$ws.Range("I137").Value = 5561564
$ws.Range("J137").Value = 1818713.9
$ws.Range("K137").Value = 16684692
$ws.Range("L137").Value = 5456141.699999999
$ws.Range("M137").Value = -16682142
$ws.Range("N137").Value = -5461241.699999999
$ws.Range("H138").Value = 2473.7727
$ws.Range("I138").Value = 1801.4375
$ws.Range("J138").Value = 4266.6665
$ws.Range("K138").Value = 5404.3125
$ws.Range("L138").Value = 12799.9995
$ws.Range("M138").Value = -264.3125
$ws.Range("N138").Value = -23079.9995
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 300
$ws.Range("I3").Value = 300
$ws.Range("K3").Value = 300
$ws.Range("M3").Value = -185
$ws.Range("H11").Value = 40005300
$ws.Range("I11").Value = 100000000
$ws.Range("J11").Value = 8833.333000000001
$ws.Range("K11").Value = 100000000
$ws.Range("L11").Value = 8833.333000000001
$ws.Range("M11").Value = -99999856
$ws.Range("N11").Value = -9121.333000000001
$ws.Range("H24").Value = 14721
$ws.Range("J24").Value = 14721
$ws.Range("L24").Value = 14721
$ws.Range("N24").Value = -15469
$ws.Range("H45").Value = 1612.75
$ws.Range("I45").Value = 1141.3334
$ws.Range("J45").Value = 1895.6
$ws.Range("K45").Value = 1141.3334
$ws.Range("L45").Value = 1895.6
$ws.Range("M45").Value = -764.3334
$ws.Range("N45").Value = -2649.6
$ws.Range("H88").Value = 18976.205
$ws.Range("I88").Value = 23972.4
$ws.Range("J88").Value = 2322.2222
$ws.Range("K88").Value = 23972.4
$ws.Range("L88").Value = 2322.2222
$ws.Range("M88").Value = -23566.4
$ws.Range("N88").Value = -3134.2222
$ws.Range("H91").Value = 18976.205
$ws.Range("I91").Value = 23972.4
$ws.Range("J91").Value = 2322.2222
$ws.Range("K91").Value = 23972.4
$ws.Range("L91").Value = 2322.2222
$ws.Range("M91").Value = -22568.4
$ws.Range("N91").Value = -5130.2222
$ws.Range("H92").Value = 26115.555
$ws.Range("J92").Value = 26115.555
$ws.Range("L92").Value = 26115.555
$ws.Range("N92").Value = -31107.555
$ws.Range("H97").Value = 854.0833
$ws.Range("I97").Value = 822.63635
$ws.Range("K97").Value = 822.63635
$ws.Range("M97").Value = -326.63635
$ws.Range("H100").Value = 14721
$ws.Range("J100").Value = 14721
$ws.Range("L100").Value = 14721
$ws.Range("N100").Value = -16885
$ws.Range("H101").Value = 37499.668
$ws.Range("J101").Value = 37499.668
$ws.Range("L101").Value = 37499.668
$ws.Range("N101").Value = -43989.668
$ws.Range("H132").Value = 42652.88
$ws.Range("I132").Value = 78724.92
$ws.Range("J132").Value = 3574.8333
$ws.Range("K132").Value = 236174.76
$ws.Range("L132").Value = 10724.4999
$ws.Range("M132").Value = -233644.76
$ws.Range("N132").Value = -15784.4999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 1074.75
$ws.Range("I11").Value = 399.5
$ws.Range("J11").Value = 1750
$ws.Range("K11").Value = 399.5
$ws.Range("L11").Value = 1750
$ws.Range("M11").Value = -259.5
$ws.Range("N11").Value = -2030
$ws.Range("H86").Value = 6436.143
$ws.Range("I86").Value = 6700.4614
$ws.Range("K86").Value = 6700.4614
$ws.Range("M86").Value = -5577.4614
$ws.Range("H89").Value = 6436.143
$ws.Range("I89").Value = 6700.4614
$ws.Range("K89").Value = 33502.307
$ws.Range("M89").Value = -27886.307
$ws.Range("H105").Value = 97094.89999999999
$ws.Range("I105").Value = 318163.34
$ws.Range("J105").Value = 2351.2856
$ws.Range("K105").Value = 318163.34
$ws.Range("L105").Value = 2351.2856
$ws.Range("M105").Value = -316416.34
$ws.Range("N105").Value = -5845.2856
$ws.Range("H107").Value = 1702.5
$ws.Range("J107").Value = 1441
$ws.Range("L107").Value = 1441
$ws.Range("N107").Value = -5281
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 548
$ws.Range("I10").Value = 548
$ws.Range("K10").Value = 548
$ws.Range("M10").Value = -409
$ws.Range("H15").Value = 3232
$ws.Range("I15").Value = 5500
$ws.Range("J15").Value = 1720
$ws.Range("K15").Value = 5500
$ws.Range("L15").Value = 1720
$ws.Range("M15").Value = -5330
$ws.Range("N15").Value = -2060
$ws.Range("H95").Value = 10579.444
$ws.Range("J95").Value = 10579.444
$ws.Range("L95").Value = 10579.444
$ws.Range("N95").Value = -16071.444
$ws.Range("H96").Value = 9649.615
$ws.Range("J96").Value = 9649.615
$ws.Range("L96").Value = 9649.615
$ws.Range("N96").Value = -15141.615
$ws.Range("H127").Value = 30000
$ws.Range("J127").Value = 30000
$ws.Range("L127").Value = 30000
$ws.Range("N127").Value = -39920
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 2724.5715
$ws.Range("I64").Value = 1961
$ws.Range("J64").Value = 3030
$ws.Range("K64").Value = 5883
$ws.Range("L64").Value = 9090
$ws.Range("M64").Value = -5613
$ws.Range("N64").Value = -9630
$ws.Range("H67").Value = 2724.5715
$ws.Range("I67").Value = 1961
$ws.Range("J67").Value = 3030
$ws.Range("K67").Value = 5883
$ws.Range("L67").Value = 9090
$ws.Range("M67").Value = -4947
$ws.Range("N67").Value = -10962
$ws.Range("H131").Value = 1119.4531
$ws.Range("I131").Value = 841.94116
$ws.Range("J131").Value = 1219.8298
$ws.Range("K131").Value = 2525.82348
$ws.Range("L131").Value = 3659.4894
$ws.Range("M131").Value = 2514.17652
$ws.Range("N131").Value = -13739.4894
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 1000
$ws.Range("I5").Value = 1000
$ws.Range("K5").Value = 1000
$ws.Range("M5").Value = -888
$ws.Range("H51").Value = 29919.4
$ws.Range("J51").Value = 29919.4
$ws.Range("L51").Value = 29919.4
$ws.Range("N51").Value = -30937.4
$ws.Range("H124").Value = 15000
$ws.Range("J124").Value = 15000
$ws.Range("L124").Value = 15000
$ws.Range("N124").Value = -24820
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2488.7778
$ws.Range("I7").Value = 2629.8
$ws.Range("J7").Value = 2312.5
$ws.Range("K7").Value = 2629.8
$ws.Range("L7").Value = 2312.5
$ws.Range("M7").Value = -2517.8
$ws.Range("N7").Value = -2536.5
$ws.Range("H82").Value = 3120.6
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 3120.6
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 3120.6
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -3842.6
$ws.Range("H85").Value = 3120.6
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 3120.6
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 3120.6
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -5616.6
$ws.Range("H126").Value = 2488.7778
$ws.Range("I126").Value = 2629.8
$ws.Range("J126").Value = 2312.5
$ws.Range("K126").Value = 7889.400000000001
$ws.Range("L126").Value = 6937.5
$ws.Range("M126").Value = -5419.400000000001
$ws.Range("N126").Value = -11877.5
$ws.Range("H132").Value = 7944222
$ws.Range("I132").Value = 3766.35
$ws.Range("J132").Value = 15162818
$ws.Range("K132").Value = 11299.05
$ws.Range("L132").Value = 45488454
$ws.Range("M132").Value = -8769.049999999999
$ws.Range("N132").Value = -45493514
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 20000
$ws.Range("J16").Value = 20000
$ws.Range("L16").Value = 20000
$ws.Range("N16").Value = -20584
$ws.Range("H132").Value = 1550.8197
$ws.Range("I132").Value = 1008.55316
$ws.Range("J132").Value = 3371.2856
$ws.Range("K132").Value = 3025.65948
$ws.Range("L132").Value = 10113.8568
$ws.Range("M132").Value = -495.6594800000003
$ws.Range("N132").Value = -15173.8568
$ws.Range("H136").Value = 16262369
$ws.Range("I136").Value = 17432272
$ws.Range("J136").Value = 10093791
$ws.Range("K136").Value = 52296816
$ws.Range("L136").Value = 30281373
$ws.Range("M136").Value = -52294266
$ws.Range("N136").Value = -30286473
